$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out the data row (row 67): remove A67, B67, C67 entirely and
# blank out D67/E67 while keeping their existing style.
$ws.Range("A67:C67").ClearContents()
$ws.Range("D67:E67").ClearContents()

# Scroll the view down so row 61 is the top visible row and select the
# range that corresponds to the now-empty row (mirrors the author's
# on-screen selection when they deleted the row's data).
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("A67:G68").Select()
